# Oakland budget data workbook - add a chunk of non-GF (non-General-Fund)
# flow rows (rows 42-53 of Sheet1) that describe additional funding
# sources flowing into various departments/agencies.
#
# Columns on Sheet1 (1-based):
#   B = index (formula, auto)
#   C = node name (source) - shared string
#   D = formula (auto) -> JSON node fragment
#   E = source index (number, lookup key into B/C)
#   F = formula (auto) -> source text
#   G = target index (number, lookup key into B/C)
#   H = formula (auto) -> target text
#   I = amount (number)
#   J = formula (auto) -> JSON link fragment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: "Kid's First Oakland Children's Fund" node name was missing;
# the E/G/I values already existed in the source workbook.
$ws.Cells.Item(42, 3).Value = "Kid's First Oakland Children's Fund"

# Row 43: Department of Health and Human Services (Federal) -> Dept of Human Services
$ws.Cells.Item(43, 3).Value = "Department of Health and Human Services (Federal)"
$ws.Cells.Item(43, 5).Value = 33
$ws.Cells.Item(43, 7).Value = 14
$ws.Cells.Item(43, 9).Value = 10756257

# Row 44: Measure B - ACTIA -> Dept of Human Services
$ws.Cells.Item(44, 3).Value = "Measure B - ACTIA"
$ws.Cells.Item(44, 5).Value = 34
$ws.Cells.Item(44, 7).Value = 14
$ws.Cells.Item(44, 9).Value = 16790670

# Row 45: Measure Q - Library Services Retention & Enhancement -> Public Works Agency
$ws.Cells.Item(45, 3).Value = "Measure Q - Library Services Retention & Enhancement"
$ws.Cells.Item(45, 5).Value = 35
$ws.Cells.Item(45, 7).Value = 15
$ws.Cells.Item(45, 9).Value = 4650750

# Row 46: Landscaping & Lighting Assessment District -> Capital Improvement Projects
$ws.Cells.Item(46, 3).Value = "Landscaping & Lighting Assessment District"
$ws.Cells.Item(46, 5).Value = 35
$ws.Cells.Item(46, 7).Value = 18
$ws.Cells.Item(46, 9).Value = 3500000

# Row 47: Development Service Fund -> Library
$ws.Cells.Item(47, 3).Value = "Development Service Fund"
$ws.Cells.Item(47, 5).Value = 36
$ws.Cells.Item(47, 7).Value = 12
$ws.Cells.Item(47, 9).Value = 14923021

# Row 48: Sewer Service Fund -> Public Works Agency
$ws.Cells.Item(48, 3).Value = "Sewer Service Fund"
$ws.Cells.Item(48, 5).Value = 37
$ws.Cells.Item(48, 7).Value = 15
$ws.Cells.Item(48, 9).Value = 14365458

# Row 49: (no new node name) Landscaping & Lighting Assessment District -> Office of Parks and Recreation
$ws.Cells.Item(49, 5).Value = 37
$ws.Cells.Item(49, 7).Value = 13
$ws.Cells.Item(49, 9).Value = 4091501

# Row 50: Development Service Fund -> Community and Economic Development Agency
$ws.Cells.Item(50, 5).Value = 38
$ws.Cells.Item(50, 7).Value = 16
$ws.Cells.Item(50, 9).Value = 21280707

# Row 51: Sewer Service Fund -> Public Works Agency
$ws.Cells.Item(51, 5).Value = 39
$ws.Cells.Item(51, 7).Value = 15
$ws.Cells.Item(51, 9).Value = 22261648

# Row 52: Sewer Service Fund -> Non Department and Port
$ws.Cells.Item(52, 5).Value = 39
$ws.Cells.Item(52, 7).Value = 17
$ws.Cells.Item(52, 9).Value = 10161330

# Row 53: Sewer Service Fund -> Capital Improvement Projects
$ws.Cells.Item(53, 5).Value = 39
$ws.Cells.Item(53, 7).Value = 18
$ws.Cells.Item(53, 9).Value = 15150000

# Columns D and J both use "look-ahead" formulas
# (IF(ISBLANK(<same-column-next-row>),...)) to decide whether to append a
# trailing comma for the JSON-ish output they build. These are shared
# formulas (cells carry only an si="" reference, no literal <f> body), and
# the engine's dependency graph for that flavor of shared formula is not
# re-evaluated just because the referenced cells' values changed later in
# the same script. Re-assigning each cell's Formula to itself forces a
# fresh recalculation against the now-populated C/E/G/I inputs.
for ($r = 42; $r -le 47; $r++) {
    $ws.Cells.Item($r, 4).Formula = $ws.Cells.Item($r, 4).Formula
}
for ($r = 43; $r -le 53; $r++) {
    $ws.Cells.Item($r, 10).Formula = $ws.Cells.Item($r, 10).Formula
}

# Match the author's final on-screen selection/view state.
$ws.Range("C49").Select()
